$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), picking up the same formatting
# (bold, bordered, centered) already used by the rest of the header row
# (e.g. H1) by copying its formats over before writing the new labels.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2 for the single data row
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
